$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 19, pushing current rows 19-31 down to 20-32.
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new weekly record.
$ws.Cells.Item(19, 1).Value = 8
$ws.Cells.Item(19, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(19, 3).Value = "Coquimbo"
$ws.Cells.Item(19, 4).Value = 45126
$ws.Cells.Item(19, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(19, 5).Value = 4
$ws.Cells.Item(19, 6).Value = 100112013
$ws.Cells.Item(19, 7).Value = "Alcachofa"
$ws.Cells.Item(19, 8).Value = "Española"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 500
$ws.Cells.Item(19, 11).Value = 14000
$ws.Cells.Item(19, 12).Value = 15000
$ws.Cells.Item(19, 13).Value = 14500
$ws.Cells.Item(19, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(19, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(19, 16).Value = 483
$ws.Cells.Item(19, 17).Value = 30
$ws.Cells.Item(19, 18).Value = "Hortaliza"
